$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F8").Value = 176
$ws.Range("F9").Value = 764
$ws.Range("F11").Value = 1049
$ws.Range("F13").Value = 793
$ws.Range("F18").Value = 1309
$ws.Range("F26").Value = 2479
$ws.Range("F27").Value = 5937
$ws.Range("F30").Value = 610
$ws.Range("F35").Value = 77
$ws.Range("F37").Value = 723
$ws.Range("F44").Value = 48
$ws.Range("F47").Value = 573
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 12
$ws.Range("F24").Value = 1714
$ws.Range("F36").Value = 100
$ws.Range("F37").Value = 176
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 781
$ws.Range("F7").Value = 251
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F7").Value = 251
$ws.Range("F8").Value = 251
$ws.Range("F12").Value = 176
$ws.Range("F14").Value = 764
$ws.Range("F18").Value = 1049
$ws.Range("F21").Value = 793
$ws.Range("F25").Value = 1309
$ws.Range("F33").Value = 2479
$ws.Range("F34").Value = 5937
$ws.Range("F36").Value = 1714
$ws.Range("F37").Value = 610
$ws.Range("F40").Value = 77
$ws.Range("F41").Value = 723
$ws.Range("F47").Value = 48
